$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 331
$ws1.Range("F7").Value = 879
$ws1.Range("F9").Value = 525
$ws1.Range("F12").Value = 1159
$ws1.Range("F21").Value = 7604
$ws1.Range("F26").Value = 2127
$ws1.Range("F27").Value = 902
$ws1.Range("F32").Value = 236
$ws1.Range("F34").Value = 1729
$ws1.Range("F36").Value = 182
$ws1.Range("F39").Value = 1225
$ws1.Range("F40").Value = 1825
$ws1.Range("F41").Value = 2147

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 7

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 331
$ws4.Range("F9").Value = 879
$ws4.Range("F11").Value = 525
$ws4.Range("F12").Value = 7
$ws4.Range("F14").Value = 1159
$ws4.Range("F24").Value = 7604
$ws4.Range("F29").Value = 2127
$ws4.Range("F30").Value = 902
$ws4.Range("F36").Value = 236
$ws4.Range("F38").Value = 1729
$ws4.Range("F40").Value = 182
$ws4.Range("F44").Value = 1225
$ws4.Range("F45").Value = 1825
$ws4.Range("F47").Value = 2147
